# Apply the "lines_states" update: the table grows from 14 data rows (A2:E15)
# to 16 data rows (A2:E17). Two new "line" entries (line7, line8) are
# effectively inserted after line6 and before the "extr" entries, which
# pushes the existing extr1..extr8 rows down by two rows; two brand new
# "extr7"/"extr8" rows are appended at the bottom with fresh values.
#
# Rather than physically inserting rows (which would introduce new style
# definitions in this runtime), we just rewrite every data cell (B:E for
# rows 2-15, plus format+populate the two new rows 16-17) to its final
# value, and extend the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for every data row: A (index), B (name), C (from_bus),
# D (to_bus), E (in_service)
$rows = @(
  @{ r = 2;  a = 0;  b = "line1"; c = 7;  d = 9;  e = $true  },
  @{ r = 3;  a = 1;  b = "line2"; c = 9;  d = 8;  e = $true  },
  @{ r = 4;  a = 2;  b = "line3"; c = 8;  d = 10; e = $true  },
  @{ r = 5;  a = 3;  b = "line4"; c = 8;  d = 11; e = $true  },
  @{ r = 6;  a = 4;  b = "line5"; c = 10; d = 5;  e = $false },
  @{ r = 7;  a = 5;  b = "line6"; c = 12; d = 8;  e = $true  },
  @{ r = 8;  a = 6;  b = "line7"; c = 14; d = 11; e = $true  },
  @{ r = 9;  a = 7;  b = "line8"; c = 16; d = 9;  e = $true  },
  @{ r = 10; a = 8;  b = "extr1"; c = 5;  d = 12; e = $false },
  @{ r = 11; a = 9;  b = "extr2"; c = 5;  d = 9;  e = $false },
  @{ r = 12; a = 10; b = "extr3"; c = 10; d = 11; e = $false },
  @{ r = 13; a = 11; b = "extr4"; c = 7;  d = 8;  e = $false },
  @{ r = 14; a = 12; b = "extr5"; c = 9;  d = 11; e = $false },
  @{ r = 15; a = 13; b = "extr6"; c = 7;  d = 11; e = $false },
  @{ r = 16; a = 14; b = "extr7"; c = 5;  d = 7;  e = $false },
  @{ r = 17; a = 15; b = "extr8"; c = 8;  d = 5;  e = $false }
)

foreach ($row in $rows) {
    $r = $row.r

    if ($r -gt 15) {
        # Brand-new row: clone formatting from the row above (column A has a
        # bold/bordered/centered style) before filling in values.
        $ws.Range("A15").Copy()
        $ws.Range("A" + $r).PasteSpecial(-4122)
    }

    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
}

$ws.Application.CutCopyMode = $false
